# Generate Report for Handoff
# Refresh the "latest handoff" timestamps for the c2f6c1fc-... file (row 7)
# across the per-locale sheets and roll the newest one up into the Overview
# sheet's "Latest HO Xliff Generate Date" column.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-25 02:41:14"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-08-25 02:41:19"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-25 02:41:19"
